# admin/course_prices.xlsx — "Price per hour changed"
#
# - Tabla1's "Precio por hora" values (I2:I4) go from 72.5 to 50.
# - The discount-table formulas (H9:H12) are generalized to reference the
#   discount percentage already typed in column G instead of a hard-coded
#   literal (G9=0.3, G10=0.2, G11=0.1, G12=0).
# - Selection moves to H9 (was M2:M4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price per hour: 72.5 -> 50 for all three course rows
$ws.Range("I2").Value = 50
$ws.Range("I3").Value = 50
$ws.Range("I4").Value = 50

# Discount-table formulas: use the % already present in column G
$ws.Range("H9").Formula  = "=J2-(L2*G9)"
$ws.Range("H10").Formula = "=L2-(L3*G10)"
$ws.Range("H11").Formula = "=M2-(M2*G11)"
$ws.Range("H12").Formula = "=K2-(K2*G12)"

$wb.Application.Calculate()

# Move the active selection to H9
$ws.Range("H9").Select() | Out-Null
